# Update cryptocurrency price/symbol data to the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'283.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'20.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'6.230"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.06190"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.586"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.559"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'1.474"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8165"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Value = "'0.1642"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.08290"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03613"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.03137"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.09140"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.709"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.001629"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.04678"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006470"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'0.001066"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.0001497"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.820"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'0.3377"
$ws.Range("D25").Style = "Normal"
$ws.Range("D40").Value = "'0.04703"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.007024"
$ws.Range("D41").Style = "Normal"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.004392"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1104"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.01138"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006372"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.9991"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.002770"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00001896"
$ws.Range("D49").Style = "Normal"
